$d = $word.ActiveDocument

# 1. Insert the two new title paragraphs before the existing (bookmark) paragraph.
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertBefore("Blade of Following`r10-Page GDD`r")

# 2. Center + bold the two new paragraphs.
$p1 = $d.Paragraphs(1)
$p1.Alignment = 1
$p1.Range.Bold = 1

$p2 = $d.Paragraphs(2)
$p2.Alignment = 1
$p2.Range.Bold = 1

# 3. Center the original (now third) paragraph -- it keeps its bookmark.
$p3 = $d.Paragraphs(3)
$p3.Alignment = 1

# 4. Drop down to raw OOXML (flat-OPC) to finish off the two changes that
#    aren't reachable through the higher-level object model in this host:
#      a) give the 3rd paragraph's mark (pPr/rPr) bold, matching the others
#      b) drop the page header entirely (sectPr headerReference + header1.xml part)
$xml = $d.WordOpenXML

# a) Bold the paragraph mark of the (now) third paragraph -- the one that still
#    carries the _GoBack bookmark and has no visible run of its own.
$bookmarkPPr = '<w:pPr><w:jc w:val="center"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/>'
$bookmarkPPrBold = '<w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/>'
$xml = $xml.Replace($bookmarkPPr, $bookmarkPPrBold)

# b) Remove the <w:headerReference .../> pointer from sectPr ...
$xml = $xml -replace '<w:headerReference[^>]*/>', ''

# ... drop the header1.xml part itself ...
$xml = $xml -replace '<pkg:part pkg:name="/word/header1\.xml"[^>]*>.*?</pkg:part>', ''

# ... and drop the now-dangling relationship that pointed at it.
$xml = $xml -replace '<Relationship[^>]*Target="header1\.xml"[^>]*/>', ''

$d.WordOpenXML = $xml
